# "finestra incidenza 7gg centrata su ultimo g"
#
# Column B holds daily new-positive counts starting on row 2.
# Columns C/D held a 7-day window CENTERED on the current day
# (days r-3 .. r+3). This script re-centers the window on the
# LAST day instead, i.e. a trailing window (days r-6 .. r):
#   C[r] = SUM(B[r-6 .. r])                     (only if all 7 days exist)
#   D[r] = C[r] / population * 100000
# Rows whose trailing window would reach above the first data
# row (row 2) are left blank, matching the existing blank rows
# at the very top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = $ws.UsedRange.Rows.Count

$population = 3463

# Cache column B values first (reading while writing C/D is fine here
# since B is never touched, but caching keeps it simple/fast).
$b = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $b[$r] = $ws.Cells.Item($r, 2).Value2
}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $windowStart = $r - 6

    if ($windowStart -lt $firstDataRow) {
        $newSum = $null
    }
    else {
        $newSum = 0
        for ($rr = $windowStart; $rr -le $r; $rr++) {
            $newSum = $newSum + $b[$rr]
        }
    }

    # Only touch cells whose effective value actually changes, so rows
    # untouched by the re-centred window keep their exact original
    # (pre-existing) numbers instead of being rewritten with freshly
    # recomputed floats that differ in the last bit.
    $oldC = $ws.Cells.Item($r, 3).Value2

    if ($newSum -eq $null) {
        $oldIsBlank = ($oldC -is [string])
        if (-not $oldIsBlank) {
            $ws.Cells.Item($r, 3).Value = ""
            $ws.Cells.Item($r, 4).Value = ""
        }
    }
    else {
        $oldIsBlank = ($oldC -is [string])
        $changed = $true
        if ((-not $oldIsBlank) -and ([Math]::Abs([double]$oldC - $newSum) -lt 0.0000001)) {
            $changed = $false
        }
        if ($changed) {
            $ws.Cells.Item($r, 3).Value = $newSum
            $ws.Cells.Item($r, 4).Value = ($newSum / $population) * 100000
        }
    }
}
